# Auto-generated edit script: update crypto price/volume table cells
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.725.82"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "2.614.43"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'601.42"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'154.46"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "2.612.28"
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("E10").Value = "  +10.91%  "
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("D12").Value = "'5.24"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("D14").Value = "'27.60"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("E15").Value = "  +3.73%  "
$ws.Range("D16").Value = "3.090.07"
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "67.662.87"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "2.612.81"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'11.15"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").Value = "'365.11"
$ws.Range("E20").Value = "  +2.77%  "
$ws.Range("D21").Value = "'7.63"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  -1.58%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'70.31"
$ws.Range("E25").Value = "  +4.37%  "
$ws.Range("D26").Value = "'9.81"
$ws.Range("E26").Value = "  -6.48%  "
$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "'0.0000104"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.745.52"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "'575.88"
$ws.Range("E29").Value = "  -3.87%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").Value = "'7.89"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  -1.71%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -3.72%  "
$ws.Range("D37").Value = "'4.93"
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D38").Value = "'158.12"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("D39").Value = "'19.36"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -2.41%  "
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'41.18"
$ws.Range("E44").Value = "  -0.71%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'16.41"
$ws.Range("D47").Value = "'156.36"
$ws.Range("E47").Value = "  +0.39%  "
$ws.Range("E48").Value = "  -7.84%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -1.99%  "
$ws.Range("E51").Value = "  +1.46%  "
